$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1868.8572
$ws.Range("J19").Value = 5291
$ws.Range("L19").Value = 5291
$ws.Range("N19").Value = -5641
$ws.Range("H41").Value = 499.58823
$ws.Range("I41").Value = 309.83334
$ws.Range("K41").Value = 309.83334
$ws.Range("M41").Value = 130.16666
$ws.Range("H103").Value = 697
$ws.Range("J103").Value = 495.66666
$ws.Range("L103").Value = 1486.99998
$ws.Range("N103").Value = -2658.99998
$ws.Range("H112").Value = 2721.2727
$ws.Range("J112").Value = 2761.875
$ws.Range("L112").Value = 8285.625
$ws.Range("N112").Value = -10501.625
$ws.Range("H132").Value = 33334388
$ws.Range("I132").Value = 41667360
$ws.Range("K132").Value = 125002080
$ws.Range("M132").Value = -124999550
$ws.Range("H137").Value = 2452.65
$ws.Range("I137").Value = 2668
$ws.Range("K137").Value = 8004
$ws.Range("M137").Value = -5454
$ws.Range("H141").Value = 3493.6155
$ws.Range("I141").Value = 3493.6155
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 10480.8465
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -5300.8465
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3278.6
$ws.Range("J63").Value = 3799.6667
$ws.Range("L63").Value = 3799.6667
$ws.Range("N63").Value = -5171.6667
$ws.Range("H66").Value = 3278.6
$ws.Range("J66").Value = 3799.6667
$ws.Range("L66").Value = 18998.3335
$ws.Range("N66").Value = -25862.3335
$ws.Range("H102").Value = 21788680
$ws.Range("I102").Value = 25001982
$ws.Range("J102").Value = 366666
$ws.Range("K102").Value = 25001982
$ws.Range("L102").Value = 366666
$ws.Range("M102").Value = -25000360
$ws.Range("N102").Value = -369910
$ws.Range("H132").Value = 5842.548
$ws.Range("I132").Value = 3507.388
$ws.Range("K132").Value = 10522.164
$ws.Range("M132").Value = -7992.164000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3274.5833
$ws.Range("I86").Value = 3359
$ws.Range("J86").Value = 3214.2856
$ws.Range("K86").Value = 3359
$ws.Range("L86").Value = 3214.2856
$ws.Range("M86").Value = -2236
$ws.Range("N86").Value = -5460.2856
$ws.Range("H89").Value = 3274.5833
$ws.Range("I89").Value = 3359
$ws.Range("J89").Value = 3214.2856
$ws.Range("K89").Value = 16795
$ws.Range("L89").Value = 16071.428
$ws.Range("M89").Value = -11179
$ws.Range("N89").Value = -27303.428
$ws.Range("H94").Value = 47621788
$ws.Range("I94").Value = 1301.2858
$ws.Range("K94").Value = 1301.2858
$ws.Range("M94").Value = -850.2858000000001
$ws.Range("H99").Value = 70176856
$ws.Range("I99").Value = 121213270
$ws.Range("J99").Value = 1772
$ws.Range("K99").Value = 121213270
$ws.Range("L99").Value = 1772
$ws.Range("M99").Value = -121211772
$ws.Range("N99").Value = -4768

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1384.6666
$ws.Range("I31").Value = 1083.7646
$ws.Range("K31").Value = 1083.7646
$ws.Range("M31").Value = -788.7646
$ws.Range("H34").Value = 1384.6666
$ws.Range("I34").Value = 1083.7646
$ws.Range("K34").Value = 1083.7646
$ws.Range("M34").Value = -881.7646
$ws.Range("H51").Value = 41998.223
$ws.Range("J51").Value = 41998.223
$ws.Range("L51").Value = 41998.223
$ws.Range("N51").Value = -43470.223
$ws.Range("H61").Value = 41998.223
$ws.Range("J61").Value = 41998.223
$ws.Range("L61").Value = 41998.223
$ws.Range("N61").Value = -42694.223
$ws.Range("H99").Value = 4463
$ws.Range("I99").Value = 2955.3333
$ws.Range("K99").Value = 2955.3333
$ws.Range("M99").Value = -1457.3333
$ws.Range("H126").Value = 4463
$ws.Range("I126").Value = 2955.3333
$ws.Range("K126").Value = 8865.999899999999
$ws.Range("M126").Value = -6395.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12821637
$ws.Range("I2").Value = 273
$ws.Range("J2").Value = 27779896
$ws.Range("K2").Value = 1638
$ws.Range("L2").Value = 166679376
$ws.Range("M2").Value = -1525
$ws.Range("N2").Value = -166679602
$ws.Range("H5").Value = 333.33334
$ws.Range("J5").Value = 400
$ws.Range("L5").Value = 1200
$ws.Range("N5").Value = -1424
$ws.Range("H80").Value = 4248.75
$ws.Range("I80").Value = 3998.5
$ws.Range("K80").Value = 11995.5
$ws.Range("M80").Value = -11059.5
$ws.Range("H83").Value = 4248.75
$ws.Range("I83").Value = 3998.5
$ws.Range("K83").Value = 35986.5
$ws.Range("M83").Value = -31306.5
$ws.Range("H134").Value = 4269.8076
$ws.Range("I134").Value = 1907.1177
$ws.Range("K134").Value = 5721.3531
$ws.Range("M134").Value = -651.3531000000003
$ws.Range("H135").Value = 333.33334
$ws.Range("J135").Value = 400
$ws.Range("L135").Value = 3600
$ws.Range("N135").Value = -8670

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2175
$ws.Range("I102").Value = 2241.4375
$ws.Range("J102").Value = 1962.4
$ws.Range("K102").Value = 2241.4375
$ws.Range("L102").Value = 1962.4
$ws.Range("M102").Value = -619.4375
$ws.Range("N102").Value = -5206.4
$ws.Range("H113").Value = 1952.3889
$ws.Range("I113").Value = 1984
$ws.Range("J113").Value = 1889.1666
$ws.Range("K113").Value = 1984
$ws.Range("L113").Value = 1889.1666
$ws.Range("M113").Value = 186
$ws.Range("N113").Value = -6229.1666
$ws.Range("H126").Value = 5636.636
$ws.Range("I126").Value = 4168.1665
$ws.Range("J126").Value = 7398.8
$ws.Range("K126").Value = 12504.4995
$ws.Range("L126").Value = 22196.4
$ws.Range("M126").Value = -10034.4995
$ws.Range("N126").Value = -27136.4
$ws.Range("H132").Value = 1635.2727
$ws.Range("J132").Value = 3499.75
$ws.Range("L132").Value = 10499.25
$ws.Range("N132").Value = -15559.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6880.346
$ws.Range("I122").Value = 4272.3184
$ws.Range("K122").Value = 12816.9552
$ws.Range("M122").Value = -10366.9552
$ws.Range("H136").Value = 2959.4666
$ws.Range("I136").Value = 2699.0833
$ws.Range("K136").Value = 8097.249899999999
$ws.Range("M136").Value = -5547.249899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2375.25
$ws.Range("I122").Value = 2332.8823
$ws.Range("K122").Value = 6998.646900000001
$ws.Range("M122").Value = -4548.646900000001
$ws.Range("H132").Value = 2273.2
$ws.Range("I132").Value = 2273.2
$ws.Range("K132").Value = 6819.599999999999
$ws.Range("M132").Value = -4289.599999999999
$ws.Range("H136").Value = 12422.692
$ws.Range("I136").Value = 12422.692
$ws.Range("K136").Value = 37268.076
$ws.Range("M136").Value = -34718.076
